# Applies the "Organizing and making a standard" edit:
#   1. Extends the "Main idea" paragraph with an extra sentence about badges.
#   2. Drops the stray _GoBack bookmark that used to sit after "StackOverflow".
#   3. Appends the new "team assignments" section (Caio/John/Ali/Malkias) plus
#      the trailing blank paragraphs, re-adding a _GoBack bookmark at the very
#      end of that new content (mirroring where Word last left the cursor).

$d = $word.ActiveDocument

# 1. " ... comments, replies, users, tags" -> add a trailing " and badges." run
$findRange = $d.Content
$findRange.Find.Execute("tags", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRange.Find.Found) {
    $findRange.Collapse(0)
    $findRange.InsertAfter(" and badges.")
}

# 2. Remove the old _GoBack bookmark (was right after "StackOverflow")
if ($d.Bookmarks.Count -ge 0) {
    try {
        $oldBookmark = $d.Bookmarks.Item("_GoBack")
        $oldBookmark.Delete()
    } catch {
    }
}

# 3. Append the new "team assignments" block at the very end of the document
$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Caio</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t>Replies component, edit replies, delete replies, Badges component, Badges increments</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>John</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t>Comments page, edit comment, delete comment, likes component</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Ali</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Main Page, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pagging</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and ordering of the main page, Create new Post, Tags</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Malkias</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">User Login, User </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sigin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, User details and edit User.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p>'
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($fragment)
